$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 3.115714333333333
$ws.Cells.Item(2, 8).Value = 9.347142999999999
$ws.Cells.Item(2, 9).Value = 0.008526392243866433
$ws.Cells.Item(2, 10).Value = 0.008526392243866435
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 34.36078833333333
$ws.Cells.Item(2, 14).Value = 103.082365
$ws.Cells.Item(2, 15).Value = 0.28490270239021
$ws.Cells.Item(2, 16).Value = 0.28490270239021
$ws.Cells.Item(2, 17).Value = 107.0584007147994
$ws.Cells.Item(2, 18).Value = 963.5256064331949
$ws.Cells.Item(2, 19).Value = 0.002429192191916473
$ws.Cells.Item(2, 20).Value = 0.002429192191916474

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 3.115714333333333
$ws.Cells.Item(3, 8).Value = 9.347142999999999
$ws.Cells.Item(3, 9).Value = 0.008526392243866433
$ws.Cells.Item(3, 10).Value = 0.008526392243866435
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 21.54461566666667
$ws.Cells.Item(3, 14).Value = 64.633847
$ws.Cells.Item(3, 15).Value = 0.17863732245739
$ws.Cells.Item(3, 16).Value = 0.1786373224573899
$ws.Cells.Item(3, 17).Value = 67.12686783879123
$ws.Cells.Item(3, 18).Value = 604.1418105491209
$ws.Cells.Item(3, 19).Value = 0.001523131880665757
$ws.Cells.Item(3, 20).Value = 0.001523131880665757

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 3.115714333333333
$ws.Cells.Item(4, 8).Value = 9.347142999999999
$ws.Cells.Item(4, 9).Value = 0.008526392243866433
$ws.Cells.Item(4, 10).Value = 0.008526392243866435
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 60.03138866666666
$ws.Cells.Item(4, 14).Value = 180.094166
$ws.Cells.Item(4, 15).Value = 0.4977506538398792
$ws.Cells.Item(4, 16).Value = 0.4977506538398792
$ws.Cells.Item(4, 17).Value = 187.0406581186375
$ws.Cells.Item(4, 18).Value = 1683.365923067738
$ws.Cells.Item(4, 19).Value = 0.004244017314279792
$ws.Cells.Item(4, 20).Value = 0.004244017314279792

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 3.115714333333333
$ws.Cells.Item(5, 8).Value = 9.347142999999999
$ws.Cells.Item(5, 9).Value = 0.008526392243866433
$ws.Cells.Item(5, 10).Value = 0.008526392243866435
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 4.668551
$ws.Cells.Item(5, 14).Value = 14.005653
$ws.Cells.Item(5, 15).Value = 0.03870932131252084
$ws.Cells.Item(5, 16).Value = 0.03870932131252084
$ws.Cells.Item(5, 17).Value = 14.54587126659766
$ws.Cells.Item(5, 18).Value = 130.912841399379
$ws.Cells.Item(5, 19).Value = 0.0003300508570044113
$ws.Cells.Item(5, 20).Value = 0.0003300508570044114

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 346.7813516666667
$ws.Cells.Item(6, 8).Value = 1040.344055
$ws.Cells.Item(6, 9).Value = 0.9489938777554333
$ws.Cells.Item(6, 10).Value = 0.9489938777554335
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 34.36078833333333
$ws.Cells.Item(6, 14).Value = 103.082365
$ws.Cells.Item(6, 15).Value = 0.28490270239021
$ws.Cells.Item(6, 16).Value = 0.28490270239021
$ws.Cells.Item(6, 17).Value = 11915.68062256556
$ws.Cells.Item(6, 18).Value = 107241.1256030901
$ws.Cells.Item(6, 19).Value = 0.2703709203242876
$ws.Cells.Item(6, 20).Value = 0.2703709203242876

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 346.7813516666667
$ws.Cells.Item(7, 8).Value = 1040.344055
$ws.Cells.Item(7, 9).Value = 0.9489938777554333
$ws.Cells.Item(7, 10).Value = 0.9489938777554335
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 21.54461566666667
$ws.Cells.Item(7, 14).Value = 64.633847
$ws.Cells.Item(7, 15).Value = 0.17863732245739
$ws.Cells.Item(7, 16).Value = 0.1786373224573899
$ws.Cells.Item(7, 17).Value = 7471.270942025511
$ws.Cells.Item(7, 18).Value = 67241.4384782296
$ws.Cells.Item(7, 19).Value = 0.1695257253506863
$ws.Cells.Item(7, 20).Value = 0.1695257253506863

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 346.7813516666667
$ws.Cells.Item(8, 8).Value = 1040.344055
$ws.Cells.Item(8, 9).Value = 0.9489938777554333
$ws.Cells.Item(8, 10).Value = 0.9489938777554335
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 60.03138866666666
$ws.Cells.Item(8, 14).Value = 180.094166
$ws.Cells.Item(8, 15).Value = 0.4977506538398792
$ws.Cells.Item(8, 16).Value = 0.4977506538398792
$ws.Cells.Item(8, 17).Value = 20817.76610425368
$ws.Cells.Item(8, 18).Value = 187359.8949382831
$ws.Cells.Item(8, 19).Value = 0.4723623231428093
$ws.Cells.Item(8, 20).Value = 0.4723623231428095

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 346.7813516666667
$ws.Cells.Item(9, 8).Value = 1040.344055
$ws.Cells.Item(9, 9).Value = 0.9489938777554333
$ws.Cells.Item(9, 10).Value = 0.9489938777554335
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 4.668551
$ws.Cells.Item(9, 14).Value = 14.005653
$ws.Cells.Item(9, 15).Value = 0.03870932131252084
$ws.Cells.Item(9, 16).Value = 0.03870932131252084
$ws.Cells.Item(9, 17).Value = 1618.966426104769
$ws.Cells.Item(9, 18).Value = 14570.69783494292
$ws.Cells.Item(9, 19).Value = 0.03673490893765019
$ws.Cells.Item(9, 20).Value = 0.0367349089376502

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 0.242641
$ws.Cells.Item(10, 8).Value = 0.727923
$ws.Cells.Item(10, 9).Value = 0.0006640057845838012
$ws.Cells.Item(10, 10).Value = 0.0006640057845838013
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 34.36078833333333
$ws.Cells.Item(10, 14).Value = 103.082365
$ws.Cells.Item(10, 15).Value = 0.28490270239021
$ws.Cells.Item(10, 16).Value = 0.28490270239021
$ws.Cells.Item(10, 17).Value = 8.337336041988333
$ws.Cells.Item(10, 18).Value = 75.03602437789499
$ws.Cells.Item(10, 19).Value = 0.0001891770424306566
$ws.Cells.Item(10, 20).Value = 0.0001891770424306566

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 0.242641
$ws.Cells.Item(11, 8).Value = 0.727923
$ws.Cells.Item(11, 9).Value = 0.0006640057845838012
$ws.Cells.Item(11, 10).Value = 0.0006640057845838013
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 21.54461566666667
$ws.Cells.Item(11, 14).Value = 64.633847
$ws.Cells.Item(11, 15).Value = 0.17863732245739
$ws.Cells.Item(11, 16).Value = 0.1786373224573899
$ws.Cells.Item(11, 17).Value = 5.227607089975667
$ws.Cells.Item(11, 18).Value = 47.048463809781
$ws.Cells.Item(11, 19).Value = 0.0001186162154542687
$ws.Cells.Item(11, 20).Value = 0.0001186162154542687

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.242641
$ws.Cells.Item(12, 8).Value = 0.727923
$ws.Cells.Item(12, 9).Value = 0.0006640057845838012
$ws.Cells.Item(12, 10).Value = 0.0006640057845838013
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 60.03138866666666
$ws.Cells.Item(12, 14).Value = 180.094166
$ws.Cells.Item(12, 15).Value = 0.4977506538398792
$ws.Cells.Item(12, 16).Value = 0.4977506538398792
$ws.Cells.Item(12, 17).Value = 14.56607617746867
$ws.Cells.Item(12, 18).Value = 131.094685597218
$ws.Cells.Item(12, 19).Value = 0.000330509313430049
$ws.Cells.Item(12, 20).Value = 0.0003305093134300491

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.242641
$ws.Cells.Item(13, 8).Value = 0.727923
$ws.Cells.Item(13, 9).Value = 0.0006640057845838012
$ws.Cells.Item(13, 10).Value = 0.0006640057845838013
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 4.668551
$ws.Cells.Item(13, 14).Value = 14.005653
$ws.Cells.Item(13, 15).Value = 0.03870932131252084
$ws.Cells.Item(13, 16).Value = 0.03870932131252084
$ws.Cells.Item(13, 17).Value = 1.132781883191
$ws.Cells.Item(13, 18).Value = 10.195036948719
$ws.Cells.Item(13, 19).Value = 0.00002570321326882686
$ws.Cells.Item(13, 20).Value = 0.00002570321326882686

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 15.28030233333333
$ws.Cells.Item(14, 8).Value = 45.840907
$ws.Cells.Item(14, 9).Value = 0.04181572421611635
$ws.Cells.Item(14, 10).Value = 0.04181572421611637
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 34.36078833333333
$ws.Cells.Item(14, 14).Value = 103.082365
$ws.Cells.Item(14, 15).Value = 0.28490270239021
$ws.Cells.Item(14, 16).Value = 0.28490270239021
$ws.Cells.Item(14, 17).Value = 525.0432341450061
$ws.Cells.Item(14, 18).Value = 4725.389107305055
$ws.Cells.Item(14, 19).Value = 0.0119134128315753
$ws.Cells.Item(14, 20).Value = 0.0119134128315753

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 15.28030233333333
$ws.Cells.Item(15, 8).Value = 45.840907
$ws.Cells.Item(15, 9).Value = 0.04181572421611635
$ws.Cells.Item(15, 10).Value = 0.04181572421611637
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 21.54461566666667
$ws.Cells.Item(15, 14).Value = 64.633847
$ws.Cells.Item(15, 15).Value = 0.17863732245739
$ws.Cells.Item(15, 16).Value = 0.1786373224573899
$ws.Cells.Item(15, 17).Value = 329.2082410421366
$ws.Cells.Item(15, 18).Value = 2962.874169379229
$ws.Cells.Item(15, 19).Value = 0.007469849010583668
$ws.Cells.Item(15, 20).Value = 0.007469849010583669

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 15.28030233333333
$ws.Cells.Item(16, 8).Value = 45.840907
$ws.Cells.Item(16, 9).Value = 0.04181572421611635
$ws.Cells.Item(16, 10).Value = 0.04181572421611637
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 60.03138866666666
$ws.Cells.Item(16, 14).Value = 180.094166
$ws.Cells.Item(16, 15).Value = 0.4977506538398792
$ws.Cells.Item(16, 16).Value = 0.4977506538398792
$ws.Cells.Item(16, 17).Value = 917.2977683165069
$ws.Cells.Item(16, 18).Value = 8255.679914848562
$ws.Cells.Item(16, 19).Value = 0.02081380406935998
$ws.Cells.Item(16, 20).Value = 0.02081380406935999

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 15.28030233333333
$ws.Cells.Item(17, 8).Value = 45.840907
$ws.Cells.Item(17, 9).Value = 0.04181572421611635
$ws.Cells.Item(17, 10).Value = 0.04181572421611637
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 4.668551
$ws.Cells.Item(17, 14).Value = 14.005653
$ws.Cells.Item(17, 15).Value = 0.03870932131252084
$ws.Cells.Item(17, 16).Value = 0.03870932131252084
$ws.Cells.Item(17, 17).Value = 71.33687073858566
$ws.Cells.Item(17, 18).Value = 642.031836647271
$ws.Cells.Item(17, 19).Value = 0.001618658304597407
$ws.Cells.Item(17, 20).Value = 0.001618658304597407
